# Apply edit: rename "P" sample IDs to "pw" naming scheme, fix mislabeled "C" group,
# and update sheet view scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rows 53..80): fix mislabeled group "B" -> "C" ---
# (Rows 53-80 were previously mislabeled with the "B" group id; they belong to "C".)
for ($row = 53; $row -le 80; $row++) {
    $ws.Cells.Item($row, 1).Value = "C"
}

# --- Column B: rename sample-id labels to the new "pw" (photometric) naming scheme ---
# Pattern: "<X> bw"        -> "BL-<X>-sw1"
# Pattern: "BL-<X>1.1-<N>" -> "BL-<X>1.1-pw-<N>"
# The four group labels are set first, followed by each group's numbered samples in order.

$ws.Cells.Item(2, 2).Value = "BL-A-sw1"
$ws.Cells.Item(81, 2).Value = "BL-D-sw1"
$ws.Cells.Item(53, 2).Value = "BL-C-sw1"
$ws.Cells.Item(27, 2).Value = "BL-B-sw1"
$ws.Cells.Item(3, 2).Value = "BL-A1.1-pw-0"
$ws.Cells.Item(4, 2).Value = "BL-A1.1-pw-1"
$ws.Cells.Item(5, 2).Value = "BL-A1.1-pw-2"
$ws.Cells.Item(6, 2).Value = "BL-A1.1-pw-3"
$ws.Cells.Item(7, 2).Value = "BL-A1.1-pw-4"
$ws.Cells.Item(8, 2).Value = "BL-A1.1-pw-5"
$ws.Cells.Item(9, 2).Value = "BL-A1.1-pw-6"
$ws.Cells.Item(10, 2).Value = "BL-A1.1-pw-7"
$ws.Cells.Item(11, 2).Value = "BL-A1.1-pw-8"
$ws.Cells.Item(12, 2).Value = "BL-A1.1-pw-9"
$ws.Cells.Item(13, 2).Value = "BL-A1.1-pw-10"
$ws.Cells.Item(14, 2).Value = "BL-A1.1-pw-11"
$ws.Cells.Item(15, 2).Value = "BL-A1.1-pw-12"
$ws.Cells.Item(16, 2).Value = "BL-A1.1-pw-13"
$ws.Cells.Item(17, 2).Value = "BL-A1.1-pw-14"
$ws.Cells.Item(18, 2).Value = "BL-A1.1-pw-15"
$ws.Cells.Item(19, 2).Value = "BL-A1.1-pw-16"
$ws.Cells.Item(20, 2).Value = "BL-A1.1-pw-17"
$ws.Cells.Item(21, 2).Value = "BL-A1.1-pw-18"
$ws.Cells.Item(22, 2).Value = "BL-A1.1-pw-19"
$ws.Cells.Item(23, 2).Value = "BL-A1.1-pw-20"
$ws.Cells.Item(24, 2).Value = "BL-A1.1-pw-21"
$ws.Cells.Item(25, 2).Value = "BL-A1.1-pw-22"
$ws.Cells.Item(26, 2).Value = "BL-A1.1-pw-23"
$ws.Cells.Item(28, 2).Value = "BL-B1.1-pw-0"
$ws.Cells.Item(29, 2).Value = "BL-B1.1-pw-1"
$ws.Cells.Item(30, 2).Value = "BL-B1.1-pw-2"
$ws.Cells.Item(31, 2).Value = "BL-B1.1-pw-3"
$ws.Cells.Item(32, 2).Value = "BL-B1.1-pw-4"
$ws.Cells.Item(33, 2).Value = "BL-B1.1-pw-5"
$ws.Cells.Item(34, 2).Value = "BL-B1.1-pw-6"
$ws.Cells.Item(35, 2).Value = "BL-B1.1-pw-7"
$ws.Cells.Item(36, 2).Value = "BL-B1.1-pw-8"
$ws.Cells.Item(37, 2).Value = "BL-B1.1-pw-9"
$ws.Cells.Item(38, 2).Value = "BL-B1.1-pw-10"
$ws.Cells.Item(39, 2).Value = "BL-B1.1-pw-11"
$ws.Cells.Item(40, 2).Value = "BL-B1.1-pw-12"
$ws.Cells.Item(41, 2).Value = "BL-B1.1-pw-13"
$ws.Cells.Item(42, 2).Value = "BL-B1.1-pw-14"
$ws.Cells.Item(43, 2).Value = "BL-B1.1-pw-15"
$ws.Cells.Item(44, 2).Value = "BL-B1.1-pw-16"
$ws.Cells.Item(45, 2).Value = "BL-B1.1-pw-17"
$ws.Cells.Item(46, 2).Value = "BL-B1.1-pw-18"
$ws.Cells.Item(47, 2).Value = "BL-B1.1-pw-19"
$ws.Cells.Item(48, 2).Value = "BL-B1.1-pw-20"
$ws.Cells.Item(49, 2).Value = "BL-B1.1-pw-21"
$ws.Cells.Item(50, 2).Value = "BL-B1.1-pw-22"
$ws.Cells.Item(51, 2).Value = "BL-B1.1-pw-23"
$ws.Cells.Item(52, 2).Value = "BL-B1.1-pw-24"
$ws.Cells.Item(54, 2).Value = "BL-C1.1-pw-0"
$ws.Cells.Item(55, 2).Value = "BL-C1.1-pw-1"
$ws.Cells.Item(56, 2).Value = "BL-C1.1-pw-2"
$ws.Cells.Item(57, 2).Value = "BL-C1.1-pw-3"
$ws.Cells.Item(58, 2).Value = "BL-C1.1-pw-4"
$ws.Cells.Item(59, 2).Value = "BL-C1.1-pw-5"
$ws.Cells.Item(60, 2).Value = "BL-C1.1-pw-6"
$ws.Cells.Item(61, 2).Value = "BL-C1.1-pw-7"
$ws.Cells.Item(62, 2).Value = "BL-C1.1-pw-8"
$ws.Cells.Item(63, 2).Value = "BL-C1.1-pw-9"
$ws.Cells.Item(64, 2).Value = "BL-C1.1-pw-10"
$ws.Cells.Item(65, 2).Value = "BL-C1.1-pw-11"
$ws.Cells.Item(66, 2).Value = "BL-C1.1-pw-12"
$ws.Cells.Item(67, 2).Value = "BL-C1.1-pw-13"
$ws.Cells.Item(68, 2).Value = "BL-C1.1-pw-14"
$ws.Cells.Item(69, 2).Value = "BL-C1.1-pw-15"
$ws.Cells.Item(70, 2).Value = "BL-C1.1-pw-16"
$ws.Cells.Item(71, 2).Value = "BL-C1.1-pw-17"
$ws.Cells.Item(72, 2).Value = "BL-C1.1-pw-18"
$ws.Cells.Item(73, 2).Value = "BL-C1.1-pw-19"
$ws.Cells.Item(74, 2).Value = "BL-C1.1-pw-20"
$ws.Cells.Item(75, 2).Value = "BL-C1.1-pw-21"
$ws.Cells.Item(76, 2).Value = "BL-C1.1-pw-22"
$ws.Cells.Item(77, 2).Value = "BL-C1.1-pw-23"
$ws.Cells.Item(78, 2).Value = "BL-C1.1-pw-24"
$ws.Cells.Item(79, 2).Value = "BL-C1.1-pw-25"
$ws.Cells.Item(80, 2).Value = "BL-C1.1-pw-26"
$ws.Cells.Item(82, 2).Value = "BL-D1.1-pw-0"
$ws.Cells.Item(83, 2).Value = "BL-D1.1-pw-1"
$ws.Cells.Item(84, 2).Value = "BL-D1.1-pw-2"
$ws.Cells.Item(85, 2).Value = "BL-D1.1-pw-3"
$ws.Cells.Item(86, 2).Value = "BL-D1.1-pw-4"
$ws.Cells.Item(87, 2).Value = "BL-D1.1-pw-5"
$ws.Cells.Item(88, 2).Value = "BL-D1.1-pw-6"
$ws.Cells.Item(89, 2).Value = "BL-D1.1-pw-7"
$ws.Cells.Item(90, 2).Value = "BL-D1.1-pw-8"
$ws.Cells.Item(91, 2).Value = "BL-D1.1-pw-9"
$ws.Cells.Item(92, 2).Value = "BL-D1.1-pw-10"
$ws.Cells.Item(93, 2).Value = "BL-D1.1-pw-11"
$ws.Cells.Item(94, 2).Value = "BL-D1.1-pw-12"
$ws.Cells.Item(95, 2).Value = "BL-D1.1-pw-13"
$ws.Cells.Item(96, 2).Value = "BL-D1.1-pw-14"
$ws.Cells.Item(97, 2).Value = "BL-D1.1-pw-15"
$ws.Cells.Item(98, 2).Value = "BL-D1.1-pw-16"
$ws.Cells.Item(99, 2).Value = "BL-D1.1-pw-17"
$ws.Cells.Item(100, 2).Value = "BL-D1.1-pw-18"
$ws.Cells.Item(101, 2).Value = "BL-D1.1-pw-19"
$ws.Cells.Item(102, 2).Value = "BL-D1.1-pw-20"
$ws.Cells.Item(103, 2).Value = "BL-D1.1-pw-21"
$ws.Cells.Item(104, 2).Value = "BL-D1.1-pw-22"
$ws.Cells.Item(105, 2).Value = "BL-D1.1-pw-23"
$ws.Cells.Item(106, 2).Value = "BL-D1.1-pw-24"
$ws.Cells.Item(107, 2).Value = "BL-D1.1-pw-25"
$ws.Cells.Item(108, 2).Value = "BL-D1.1-pw-26"

# --- Sheet view: scroll position and selection, matching the saved view state ---
$ws.Application.ActiveWindow.ScrollRow = 57
$ws.Range("G118:H122").Select()
